# Backup QR Scanner data - rename the sheet and remove the latest
# manual-entry log row, matching the refreshed scanner export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the worksheet tab from "Scanner" to "Pediatrics"
$ws.Name = "Pediatrics"

# 2) Remove the third data row (student 190101 / Manual entry),
#    shifting the rest of the sheet up. This also shrinks the
#    worksheet's used range/dimension from A1:F3 down to A1:F2.
$ws.Rows.Item(3).Delete()

# 3) Re-flag the "number stored as text" warning as ignored for the
#    (now smaller) data range, matching the updated ignoredErrors
#    sqref for A1:F2.
$ws.Range("A1:F2").Errors.Item(6).Ignore = $true
